$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column B for the "Sample ID" field
$ws.Columns("B:B").Insert()

# Header
$ws.Cells.Item(1, 2).Value = "Sample ID"

# Fill numeric Sample ID values for rows 2-27 (= A value - 3000)
for ($r = 2; $r -le 27; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 2).Value = $a - 3000
}

# Rows 28-31 get zero-padded text Sample IDs, right-aligned, text number format
$first = $ws.Cells.Item(28, 2)
$first.NumberFormat = "@"
$first.HorizontalAlignment = -4152
$first.Value = "057"

$first.Copy()
$ws.Range("B29:B31").PasteSpecial(-4122)
$ws.Cells.Item(29, 2).Value = "058"
$ws.Cells.Item(30, 2).Value = "059"
$ws.Cells.Item(31, 2).Value = "060"

$ws.Columns("B:B").AutoFit()

$ws.Range("G15").Select()
